$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold + bordered + centered) onto the two
# new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-23 for the new columns I (I0) and J (IF)
$data = @{
    2  = @(1, 3)
    3  = @(1, 5)
    4  = @(1, 6)
    5  = @(1, 6)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 5)
    9  = @(1, 4)
    10 = @(1, 3)
    11 = @(8, 9)
    12 = @(8, 8)
    13 = @(8, 9)
    14 = @(4, 6)
    15 = @(2, 2)
    16 = @(4, 6)
    17 = @(7, 8)
    18 = @(8, 8)
    19 = @(6, 6)
    20 = @(5, 6)
    21 = @(6, 6)
    22 = @(3, 4)
    23 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
